# Applies the "Horarios actualizados Linea 141 - 201" update.
# Sheet 1 (LP1912): header refresh + rows 60-90 rewritten (new scrape pass
#   07:46:15 interleaved with the existing rows, total filas 64 -> 85).
# Sheet 2 (LP1912-215): header refresh + 2 new rows appended (13 -> 15).
# Sheet 3 (6203-6173): header refresh + rows 17-25 rewritten (16 -> 20).

$wb = $excel.ActiveWorkbook

$lastUpdate = '07:46:15'

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item('LP1912')
$ws1.Range('A2').Value2 = 'Última actualización: ' + $lastUpdate
$ws1.Range('A3').Value2 = 'Total filas: 85'

$data1 = New-Object 'object[,]' 31,5
$data1[0,0] = '07:46:15'
$data1[0,1] = '07:50'
$data1[0,2] = '14_ABASTO'
$data1[0,3] = 4
$data1[0,4] = 'LP1912'
$data1[1,0] = '06:23:52'
$data1[1,1] = '07:51'
$data1[1,2] = '215D_EL PATO'
$data1[1,3] = 88
$data1[1,4] = 'LP1912'
$data1[2,0] = '07:46:15'
$data1[2,1] = '07:55'
$data1[2,2] = '10_OLMOS'
$data1[2,3] = 9
$data1[2,4] = 'LP1912'
$data1[3,0] = '07:46:15'
$data1[3,1] = '07:56'
$data1[3,2] = '16_SANTA ANA'
$data1[3,3] = 10
$data1[3,4] = 'LP1912'
$data1[4,0] = '06:52:23'
$data1[4,1] = '08:03'
$data1[4,2] = '23_HERNANDEZ'
$data1[4,3] = 71
$data1[4,4] = 'LP1912'
$data1[5,0] = '06:23:52'
$data1[5,1] = '08:05'
$data1[5,2] = '23_HERNANDEZ'
$data1[5,3] = 102
$data1[5,4] = 'LP1912'
$data1[6,0] = '07:46:15'
$data1[6,1] = '08:09'
$data1[6,2] = '11_ETCHEVERRY'
$data1[6,3] = 23
$data1[6,4] = 'LP1912'
$data1[7,0] = '06:23:52'
$data1[7,1] = '08:12'
$data1[7,2] = '15_ABASTO'
$data1[7,3] = 109
$data1[7,4] = 'LP1912'
$data1[8,0] = '06:23:52'
$data1[8,1] = '08:20'
$data1[8,2] = '26_HERNANDEZ'
$data1[8,3] = 117
$data1[8,4] = 'LP1912'
$data1[9,0] = '06:52:23'
$data1[9,1] = '08:21'
$data1[9,2] = '26_HERNANDEZ'
$data1[9,3] = 89
$data1[9,4] = 'LP1912'
$data1[10,0] = '06:23:52'
$data1[10,1] = '08:22'
$data1[10,2] = '16_P MOR-SANTA ANA'
$data1[10,3] = 119
$data1[10,4] = 'LP1912'
$data1[11,0] = '06:52:23'
$data1[11,1] = '08:23'
$data1[11,2] = '215B_EL PATO'
$data1[11,3] = 91
$data1[11,4] = 'LP1912'
$data1[12,0] = '07:46:15'
$data1[12,1] = '08:23'
$data1[12,2] = '16_P MOR-SANTA ANA'
$data1[12,3] = 37
$data1[12,4] = 'LP1912'
$data1[13,0] = '06:52:23'
$data1[13,1] = '08:27'
$data1[13,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[13,3] = 95
$data1[13,4] = 'LP1912'
$data1[14,0] = '07:46:15'
$data1[14,1] = '08:33'
$data1[14,2] = '10_OLMOS'
$data1[14,3] = 47
$data1[14,4] = 'LP1912'
$data1[15,0] = '07:46:15'
$data1[15,1] = '08:33'
$data1[15,2] = '23_HERNANDEZ'
$data1[15,3] = 47
$data1[15,4] = 'LP1912'
$data1[16,0] = '07:46:15'
$data1[16,1] = '08:34'
$data1[16,2] = '26_HERNANDEZ'
$data1[16,3] = 48
$data1[16,4] = 'LP1912'
$data1[17,0] = '06:52:23'
$data1[17,1] = '08:42'
$data1[17,2] = '81_EL PELIGRO'
$data1[17,3] = 110
$data1[17,4] = 'LP1912'
$data1[18,0] = '07:46:15'
$data1[18,1] = '08:44'
$data1[18,2] = '14_ABASTO'
$data1[18,3] = 58
$data1[18,4] = 'LP1912'
$data1[19,0] = '07:46:15'
$data1[19,1] = '08:54'
$data1[19,2] = '17_ROMERO'
$data1[19,3] = 68
$data1[19,4] = 'LP1912'
$data1[20,0] = '07:46:15'
$data1[20,1] = '09:02'
$data1[20,2] = '215A_EL PATO'
$data1[20,3] = 76
$data1[20,4] = 'LP1912'
$data1[21,0] = '07:46:15'
$data1[21,1] = '09:04'
$data1[21,2] = '11_ETCHEVERRY'
$data1[21,3] = 78
$data1[21,4] = 'LP1912'
$data1[22,0] = '07:46:15'
$data1[22,1] = '09:11'
$data1[22,2] = '16_P MOR-SANTA ANA'
$data1[22,3] = 85
$data1[22,4] = 'LP1912'
$data1[23,0] = '07:46:15'
$data1[23,1] = '09:17'
$data1[23,2] = '27_EL RETIRO'
$data1[23,3] = 91
$data1[23,4] = 'LP1912'
$data1[24,0] = '07:46:15'
$data1[24,1] = '09:21'
$data1[24,2] = '26_HERNANDEZ'
$data1[24,3] = 95
$data1[24,4] = 'LP1912'
$data1[25,0] = '07:46:15'
$data1[25,1] = '09:22'
$data1[25,2] = '16_SANTA ANA'
$data1[25,3] = 96
$data1[25,4] = 'LP1912'
$data1[26,0] = '07:46:15'
$data1[26,1] = '09:23'
$data1[26,2] = '17_ROMERO'
$data1[26,3] = 97
$data1[26,4] = 'LP1912'
$data1[27,0] = '07:46:15'
$data1[27,1] = '09:24'
$data1[27,2] = '11_ETCHEVERRY'
$data1[27,3] = 98
$data1[27,4] = 'LP1912'
$data1[28,0] = '07:46:15'
$data1[28,1] = '09:32'
$data1[28,2] = '15_ABASTO'
$data1[28,3] = 106
$data1[28,4] = 'LP1912'
$data1[29,0] = '07:46:15'
$data1[29,1] = '09:33'
$data1[29,2] = '10_OLMOS'
$data1[29,3] = 107
$data1[29,4] = 'LP1912'
$data1[30,0] = '07:46:15'
$data1[30,1] = '09:42'
$data1[30,2] = '215C_EL PATO'
$data1[30,3] = 116
$data1[30,4] = 'LP1912'
$ws1.Range('A60:E90').Value2 = $data1

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws2.Range('A2').Value2 = 'Última actualización: ' + $lastUpdate
$ws2.Range('A3').Value2 = 'Total filas: 15'

$data2 = New-Object 'object[,]' 2,5
$data2[0,0] = '07:46:15'
$data2[0,1] = '09:02'
$data2[0,2] = '215A_EL PATO'
$data2[0,3] = 76
$data2[0,4] = 'LP1912'
$data2[1,0] = '07:46:15'
$data2[1,1] = '09:42'
$data2[1,2] = '215C_EL PATO'
$data2[1,3] = 116
$data2[1,4] = 'LP1912'
$ws2.Range('A19:E20').Value2 = $data2

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item('6203-6173')
$ws3.Range('A2').Value2 = 'Última actualización: ' + $lastUpdate
$ws3.Range('A3').Value2 = 'Total filas: 20'

$data3 = New-Object 'object[,]' 9,5
$data3[0,0] = '07:46:15'
$data3[0,1] = '07:51'
$data3[0,2] = '215A_LA PLATA'
$data3[0,3] = 5
$data3[0,4] = 'L6173'
$data3[1,0] = '06:23:52'
$data3[1,1] = '08:06'
$data3[1,2] = '215C_LA PLATA'
$data3[1,3] = 103
$data3[1,4] = 'L6203'
$data3[2,0] = '07:46:15'
$data3[2,1] = '08:09'
$data3[2,2] = '215C_LA PLATA'
$data3[2,3] = 23
$data3[2,4] = 'L6203'
$data3[3,0] = '07:16:53'
$data3[3,1] = '08:10'
$data3[3,2] = '215C_LA PLATA'
$data3[3,3] = 54
$data3[3,4] = 'L6203'
$data3[4,0] = '06:52:23'
$data3[4,1] = '08:11'
$data3[4,2] = '215C_LA PLATA'
$data3[4,3] = 79
$data3[4,4] = 'L6203'
$data3[5,0] = '06:52:23'
$data3[5,1] = '08:40'
$data3[5,2] = '215A_LA PLATA'
$data3[5,3] = 108
$data3[5,4] = 'L6173'
$data3[6,0] = '07:46:15'
$data3[6,1] = '08:45'
$data3[6,2] = '215A_LA PLATA'
$data3[6,3] = 59
$data3[6,4] = 'L6173'
$data3[7,0] = '07:16:53'
$data3[7,1] = '09:08'
$data3[7,2] = '215D_LA PLATA'
$data3[7,3] = 112
$data3[7,4] = 'LP1912'
$data3[8,0] = '07:46:15'
$data3[8,1] = '09:09'
$data3[8,2] = '215D_LA PLATA'
$data3[8,3] = 83
$data3[8,4] = 'L6203'
$ws3.Range('A17:E25').Value2 = $data3
